$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the FvFm (column F) values for the T_0 baseline readings (rows 2-6)
$ws.Range("F2").Value = 0.486
$ws.Range("F3").Value = 0.462
$ws.Range("F4").Value = 0.455
$ws.Range("F5").Value = 0.475
$ws.Range("F6").Value = 0.468

# Update the active selection to F7, matching the saved view state
$ws.Range("F7").Select()
